$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 422.925
$ws.Range("J17").Value = 422.925
$ws.Range("L17").Value = 1268.775
$ws.Range("N17").Value = -1604.775
$ws.Range("H33").Value = 544
$ws.Range("I33").Value = 550.4
$ws.Range("K33").Value = 550.4
$ws.Range("M33").Value = -321.4
$ws.Range("H112").Value = 2171.465
$ws.Range("I112").Value = 809.875
$ws.Range("J112").Value = 2482.6858
$ws.Range("K112").Value = 2429.625
$ws.Range("L112").Value = 7448.057400000001
$ws.Range("M112").Value = -1321.625
$ws.Range("N112").Value = -9664.057400000002
$ws.Range("H129").Value = 894.5599999999999
$ws.Range("J129").Value = 943.6957
$ws.Range("L129").Value = 2831.0871
$ws.Range("N129").Value = -12831.0871
$ws.Range("H132").Value = 18529976
$ws.Range("I132").Value = 27788388
$ws.Range("K132").Value = 83365164
$ws.Range("M132").Value = -83362634
$ws.Range("H135").Value = 35715028
$ws.Range("I135").Value = 351.85715
$ws.Range("J135").Value = 142859060
$ws.Range("K135").Value = 3166.71435
$ws.Range("L135").Value = 1285731540
$ws.Range("M135").Value = -631.7143499999997
$ws.Range("N135").Value = -1285736610
$ws.Range("H137").Value = 1422.619
$ws.Range("I137").Value = 934
$ws.Range("J137").Value = 1789.0834
$ws.Range("K137").Value = 2802
$ws.Range("L137").Value = 5367.2502
$ws.Range("M137").Value = -252
$ws.Range("N137").Value = -10467.2502
$ws.Range("H138").Value = 1419.76
$ws.Range("I138").Value = 603.881
$ws.Range("J138").Value = 2010.569
$ws.Range("K138").Value = 1811.643
$ws.Range("L138").Value = 6031.707
$ws.Range("M138").Value = 3328.357
$ws.Range("N138").Value = -16311.707

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1228.84
$ws.Range("I61").Value = 1176.35
$ws.Range("K61").Value = 1176.35
$ws.Range("M61").Value = -964.3499999999999
$ws.Range("H97").Value = 383.91666
$ws.Range("I97").Value = 391.54544
$ws.Range("J97").Value = 300
$ws.Range("K97").Value = 391.54544
$ws.Range("L97").Value = 300
$ws.Range("M97").Value = 104.45456
$ws.Range("N97").Value = -1292
$ws.Range("H122").Value = 522
$ws.Range("I122").Value = 522
$ws.Range("K122").Value = 1566
$ws.Range("M122").Value = 884
$ws.Range("H132").Value = 1729.3793
$ws.Range("I132").Value = 1358.8695
$ws.Range("J132").Value = 3149.6667
$ws.Range("K132").Value = 4076.6085
$ws.Range("L132").Value = 9449.000100000001
$ws.Range("M132").Value = -1546.6085
$ws.Range("N132").Value = -14509.0001
$ws.Range("H136").Value = 1228.84
$ws.Range("I136").Value = 1176.35
$ws.Range("K136").Value = 3529.05
$ws.Range("M136").Value = -979.0499999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1450.1333
$ws.Range("I20").Value = 1236.5714
$ws.Range("J20").Value = 1637
$ws.Range("K20").Value = 1236.5714
$ws.Range("L20").Value = 1637
$ws.Range("M20").Value = -989.5714
$ws.Range("N20").Value = -2131
$ws.Range("H99").Value = 27778762
$ws.Range("I99").Value = 45455480
$ws.Range("J99").Value = 1061.7142
$ws.Range("K99").Value = 45455480
$ws.Range("L99").Value = 1061.7142
$ws.Range("M99").Value = -45453982
$ws.Range("N99").Value = -4057.7142
$ws.Range("H105").Value = 142860450
$ws.Range("I105").Value = 142860450
$ws.Range("K105").Value = 142860450
$ws.Range("M105").Value = -142858703
$ws.Range("H122").Value = 40666.668
$ws.Range("J122").Value = 40666.668
$ws.Range("L122").Value = 40666.668
$ws.Range("N122").Value = -50466.668
$ws.Range("H124").Value = 40780
$ws.Range("J124").Value = 40780
$ws.Range("L124").Value = 40780
$ws.Range("N124").Value = -50600
$ws.Range("H134").Value = 3060.7144
$ws.Range("I134").Value = 870.02325
$ws.Range("K134").Value = 2610.06975
$ws.Range("M134").Value = -75.06974999999966

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49999
$ws.Range("J20").Value = 49999
$ws.Range("L20").Value = 49999
$ws.Range("N20").Value = -50471
$ws.Range("H30").Value = 49999
$ws.Range("J30").Value = 49999
$ws.Range("L30").Value = 49999
$ws.Range("N30").Value = -50181
$ws.Range("H31").Value = 1180.7354
$ws.Range("I31").Value = 1126.8036
$ws.Range("J31").Value = 1432.4166
$ws.Range("K31").Value = 1126.8036
$ws.Range("L31").Value = 1432.4166
$ws.Range("M31").Value = -831.8036
$ws.Range("N31").Value = -2022.4166
$ws.Range("H34").Value = 1180.7354
$ws.Range("I34").Value = 1126.8036
$ws.Range("J34").Value = 1432.4166
$ws.Range("K34").Value = 1126.8036
$ws.Range("L34").Value = 1432.4166
$ws.Range("M34").Value = -924.8036
$ws.Range("N34").Value = -1836.4166
$ws.Range("H58").Value = 1342.5294
$ws.Range("I58").Value = 1321.3334
$ws.Range("K58").Value = 1321.3334
$ws.Range("M58").Value = -1118.3334
$ws.Range("H105").Value = 884
$ws.Range("I105").Value = 806.6667
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 806.6667
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 940.3333
$ws.Range("N105").Value = -4494
$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959
$ws.Range("H129").Value = 48749.625
$ws.Range("J129").Value = 48749.625
$ws.Range("L129").Value = 48749.625
$ws.Range("N129").Value = -58749.625
$ws.Range("H132").Value = 2045.2
$ws.Range("I132").Value = 742.2222
$ws.Range("K132").Value = 2226.6666
$ws.Range("M132").Value = 303.3334
$ws.Range("H134").Value = 1288.3529
$ws.Range("I134").Value = 996.8889
$ws.Range("K134").Value = 2990.6667
$ws.Range("M134").Value = -455.6667000000002
$ws.Range("H136").Value = 1342.5294
$ws.Range("I136").Value = 1321.3334
$ws.Range("K136").Value = 3964.0002
$ws.Range("M136").Value = -1414.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5000
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -16872
$ws.Range("H83").Value = 5000
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -54360
$ws.Range("H131").Value = 14288081
$ws.Range("J131").Value = 2645.1968
$ws.Range("L131").Value = 7935.590400000001
$ws.Range("N131").Value = -18015.5904

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 64288772
$ws.Range("I70").Value = 50003360
$ws.Range("K70").Value = 50003360
$ws.Range("M70").Value = -50003090
$ws.Range("H73").Value = 64288772
$ws.Range("I73").Value = 50003360
$ws.Range("K73").Value = 50003360
$ws.Range("M73").Value = -50002424
$ws.Range("H132").Value = 2080.95
$ws.Range("I132").Value = 1537.9048
$ws.Range("J132").Value = 2681.158
$ws.Range("K132").Value = 4613.7144
$ws.Range("L132").Value = 8043.474
$ws.Range("M132").Value = -2083.7144
$ws.Range("N132").Value = -13103.474

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 845
$ws.Range("I16").Value = 845
$ws.Range("K16").Value = 845
$ws.Range("M16").Value = -675
$ws.Range("H100").Value = 1652.5
$ws.Range("I100").Value = 805
$ws.Range("K100").Value = 805
$ws.Range("M100").Value = -264
$ws.Range("H123").Value = 52000
$ws.Range("J123").Value = 52000
$ws.Range("L123").Value = 52000
$ws.Range("N123").Value = -61800
$ws.Range("H132").Value = 30333.057
$ws.Range("I132").Value = 1097.579
$ws.Range("J132").Value = 65050.188
$ws.Range("K132").Value = 3292.737
$ws.Range("L132").Value = 195150.564
$ws.Range("M132").Value = -762.7370000000001
$ws.Range("N132").Value = -200210.564
$ws.Range("H136").Value = 1991.4
$ws.Range("I136").Value = 2834.6667
$ws.Range("K136").Value = 8504.000100000001
$ws.Range("M136").Value = -5954.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1414.4375
$ws.Range("I132").Value = 1167.7241
$ws.Range("K132").Value = 3503.1723
$ws.Range("M132").Value = -973.1722999999997
$ws.Range("H136").Value = 943.4
$ws.Range("I136").Value = 631.5
$ws.Range("K136").Value = 1894.5
$ws.Range("M136").Value = 655.5
